$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.738.00"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.873.60"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4920"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2890"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06565"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").Value = "1.879.50"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07157"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6633"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.794"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "29.753.66"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007788"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").Value = "2.121.76"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.719"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.076"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.534"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.917"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.146"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08649"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.889"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05050"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7012"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.100"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.672"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.690"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.196"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9319"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01630"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.041"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9959"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4131"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.423"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1252"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.197"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.331"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "
